$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force Text format so numeric-looking strings (e.g. "0.9984", "1.0000")
# are preserved exactly as literal text instead of being parsed into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '26.337.05'
$ws.Range("E2").Value = '  +3.68%  '
$ws.Range("D3").Value = '1.717.79'
$ws.Range("E3").Value = '  +3.17%  '
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '239.40'
$ws.Range("E5").Value = '  +1.33%  '
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = '0.4709'
$ws.Range("E7").Value = '  -2.01%  '
$ws.Range("D8").Value = '0.2638'
$ws.Range("E8").Value = '  +1.13%  '
$ws.Range("D9").Value = '0.06223'
$ws.Range("E9").Value = '  +1.25%  '
$ws.Range("D10").Value = '1.710.70'
$ws.Range("E10").Value = '  +2.84%  '
$ws.Range("D11").Value = '0.07075'
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("D12").Value = '15.23'
$ws.Range("E12").Value = '  +3.41%  '
$ws.Range("D13").Value = '4.419'
$ws.Range("E13").Value = '  +0.94%  '
$ws.Range("D14").Value = '0.5899'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("D15").Value = '76.33'
$ws.Range("E15").Value = '  +2.65%  '
$ws.Range("D16").Value = '1.0000'
$ws.Range("E16").Value = '  +0.11%  '
$ws.Range("D17").Value = '0.9992'
$ws.Range("E17").Value = '  +0.09%  '
$ws.Range("D18").Value = '26.318.41'
$ws.Range("E18").Value = '  +3.63%  '
$ws.Range("D19").Value = '0.000006815'
$ws.Range("E19").Value = '  +1.80%  '
$ws.Range("D21").Value = '1.929.65'
$ws.Range("E21").Value = '  +3.25%  '
$ws.Range("D22").Value = '4.555'
$ws.Range("E22").Value = '  +3.52%  '
$ws.Range("D23").Value = '8.820'
$ws.Range("E23").Value = '  +2.24%  '
$ws.Range("D24").Value = '5.348'
$ws.Range("E24").Value = '  +0.38%  '
$ws.Range("D25").Value = '135.39'
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").Value = '15.19'
$ws.Range("E26").Value = '  +0.49%  '
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("D28").Value = '1.762'
$ws.Range("E28").Value = '  +4.14%  '
$ws.Range("D29").Value = '106.88'
$ws.Range("E29").Value = '  +2.11%  '
$ws.Range("D30").Value = '4.048'
$ws.Range("E30").Value = '  +1.53%  '
$ws.Range("D31").Value = '3.688'
$ws.Range("E31").Value = '  +2.10%  '
$ws.Range("D32").Value = '0.07717'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("D33").Value = '0.04418'
$ws.Range("E33").Value = '  +0.86%  '
$ws.Range("D34").Value = '2.612'
$ws.Range("E34").Value = '  +0.62%  '
$ws.Range("D35").Value = '0.6227'
$ws.Range("E35").Value = '  +3.12%  '
$ws.Range("D36").Value = '0.9710'
$ws.Range("E36").Value = '  +2.97%  '
$ws.Range("D37").Value = '0.9220'
$ws.Range("E37").Value = '  +7.78%  '
$ws.Range("D38").Value = '113.32'
$ws.Range("E38").Value = '  +14.91%  '
$ws.Range("D39").Value = '2.414'
$ws.Range("E39").Value = '  -7.81%  '
$ws.Range("D40").Value = '0.9997'
$ws.Range("E40").Value = '  +0.03%  '
$ws.Range("D41").Value = '1.905'
$ws.Range("E41").Value = '  +4.23%  '
$ws.Range("D42").Value = '0.01466'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '5.264'
$ws.Range("E43").Value = '  +11.98%  '
$ws.Range("D44").Value = '0.3810'
$ws.Range("D45").Value = '0.1150'
$ws.Range("E45").Value = '  +3.20%  '
$ws.Range("D46").Value = '6.243'
$ws.Range("E46").Value = '  +0.56%  '
$ws.Range("D47").Value = '0.05288'
$ws.Range("E47").Value = '  +0.79%  '
$ws.Range("D48").Value = '30.52'
$ws.Range("E48").Value = '  +3.53%  '
$ws.Range("D49").Value = '7.632'
$ws.Range("E49").Value = '  +4.88%  '
$ws.Range("D50").Value = '1.221'
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("D51").Value = '0.3380'
$ws.Range("E51").Value = '  +1.17%  '

# Remove the temporary formatting so the cells end up with no explicit style,
# matching the original look (text type is preserved on save).
$ws.Range("D2:E51").ClearFormats()
